$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 316.5
$ws.Range("I33").Value = 324.5238
$ws.Range("K33").Value = 324.5238
$ws.Range("M33").Value = -95.52379999999999
$ws.Range("H76").Value = 5573.364
$ws.Range("I76").Value = 6150.375
$ws.Range("J76").Value = 4034.6667
$ws.Range("K76").Value = 6150.375
$ws.Range("L76").Value = 4034.6667
$ws.Range("M76").Value = -5835.375
$ws.Range("N76").Value = -4664.6667
$ws.Range("H79").Value = 5573.364
$ws.Range("I79").Value = 6150.375
$ws.Range("J79").Value = 4034.6667
$ws.Range("K79").Value = 6150.375
$ws.Range("L79").Value = 4034.6667
$ws.Range("M79").Value = -5058.375
$ws.Range("N79").Value = -6218.6667
$ws.Range("H86").Value = 2149.1875
$ws.Range("I86").Value = 2274.4167
$ws.Range("J86").Value = 1773.5
$ws.Range("K86").Value = 2274.4167
$ws.Range("L86").Value = 1773.5
$ws.Range("M86").Value = -1151.4167
$ws.Range("N86").Value = -4019.5
$ws.Range("H88").Value = 1952.8462
$ws.Range("I88").Value = 650
$ws.Range("K88").Value = 650
$ws.Range("M88").Value = -244
$ws.Range("H89").Value = 2149.1875
$ws.Range("I89").Value = 2274.4167
$ws.Range("J89").Value = 1773.5
$ws.Range("K89").Value = 11372.0835
$ws.Range("L89").Value = 8867.5
$ws.Range("M89").Value = -5756.083500000001
$ws.Range("N89").Value = -20099.5
$ws.Range("H91").Value = 1952.8462
$ws.Range("I91").Value = 650
$ws.Range("K91").Value = 650
$ws.Range("M91").Value = 754
$ws.Range("H107").Value = 443.95456
$ws.Range("J107").Value = 506.30768
$ws.Range("L107").Value = 506.30768
$ws.Range("N107").Value = -4346.30768
$ws.Range("H127").Value = 1279.2106
$ws.Range("I127").Value = 480.5
$ws.Range("J127").Value = 2166.6667
$ws.Range("K127").Value = 1441.5
$ws.Range("L127").Value = 6500.000100000001
$ws.Range("M127").Value = 3518.5
$ws.Range("N127").Value = -16420.0001
$ws.Range("H129").Value = 1097.2593
$ws.Range("J129").Value = 1145.22
$ws.Range("L129").Value = 3435.66
$ws.Range("N129").Value = -13435.66
$ws.Range("H138").Value = 2686.2546
$ws.Range("I138").Value = 1723.2273
$ws.Range("J138").Value = 3328.2727
$ws.Range("K138").Value = 5169.6819
$ws.Range("L138").Value = 9984.8181
$ws.Range("M138").Value = -29.68189999999959
$ws.Range("N138").Value = -20264.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1703.3
$ws.Range("I2").Value = 1133.2858
$ws.Range("J2").Value = 3033.3333
$ws.Range("K2").Value = 1133.2858
$ws.Range("L2").Value = 3033.3333
$ws.Range("M2").Value = -1020.2858
$ws.Range("N2").Value = -3259.3333
$ws.Range("H74").Value = 1217.2667
$ws.Range("I74").Value = 912.7143
$ws.Range("J74").Value = 1483.75
$ws.Range("K74").Value = 912.7143
$ws.Range("L74").Value = 1483.75
$ws.Range("M74").Value = -38.71429999999998
$ws.Range("N74").Value = -3231.75
$ws.Range("H77").Value = 1217.2667
$ws.Range("I77").Value = 912.7143
$ws.Range("J77").Value = 1483.75
$ws.Range("K77").Value = 4563.5715
$ws.Range("L77").Value = 7418.75
$ws.Range("M77").Value = -195.5715
$ws.Range("N77").Value = -16154.75
$ws.Range("H95").Value = 25208
$ws.Range("J95").Value = 25208
$ws.Range("L95").Value = 25208
$ws.Range("N95").Value = -30700
$ws.Range("H116").Value = 1703.3
$ws.Range("I116").Value = 1133.2858
$ws.Range("J116").Value = 3033.3333
$ws.Range("K116").Value = 1133.2858
$ws.Range("L116").Value = 3033.3333
$ws.Range("M116").Value = 1160.7142
$ws.Range("N116").Value = -7621.3333
$ws.Range("H122").Value = 3246.5715
$ws.Range("I122").Value = 2864.2666
$ws.Range("J122").Value = 4202.3335
$ws.Range("K122").Value = 8592.799800000001
$ws.Range("L122").Value = 12607.0005
$ws.Range("M122").Value = -6142.799800000001
$ws.Range("N122").Value = -17507.0005
$ws.Range("H123").Value = 24172.428
$ws.Range("J123").Value = 24172.428
$ws.Range("L123").Value = 24172.428
$ws.Range("N123").Value = -33972.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1964.06
$ws.Range("I134").Value = 1558.7949
$ws.Range("J134").Value = 3400.9092
$ws.Range("K134").Value = 4676.384700000001
$ws.Range("L134").Value = 10202.7276
$ws.Range("M134").Value = -2141.384700000001
$ws.Range("N134").Value = -15272.7276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 187.5
$ws.Range("I19").Value = 187.5
$ws.Range("K19").Value = 187.5
$ws.Range("M19").Value = -17.5
$ws.Range("H24").Value = 187.5
$ws.Range("I24").Value = 187.5
$ws.Range("K24").Value = 187.5
$ws.Range("M24").Value = -17.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 630.05554
$ws.Range("I117").Value = 313.6
$ws.Range("J117").Value = 751.7692
$ws.Range("K117").Value = 940.8000000000001
$ws.Range("L117").Value = 2255.3076
$ws.Range("M117").Value = 2501.2
$ws.Range("N117").Value = -9139.3076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H62").Value = 17000
$ws.Range("I62").Value = 17000
$ws.Range("K62").Value = 17000
$ws.Range("M62").Value = -16314
$ws.Range("H65").Value = 17000
$ws.Range("I65").Value = 17000
$ws.Range("K65").Value = 51000
$ws.Range("M65").Value = -47568
$ws.Range("H80").Value = 3160
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3320
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3320
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -5316
$ws.Range("H83").Value = 3160
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3320
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 16600
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -26584
$ws.Range("H92").Value = 7257.4287
$ws.Range("J92").Value = 7257.4287
$ws.Range("L92").Value = 7257.4287
$ws.Range("N92").Value = -11001.4287
$ws.Range("H109").Value = 9311.125
$ws.Range("J109").Value = 9311.125
$ws.Range("L109").Value = 9311.125
$ws.Range("N109").Value = -11391.125
$ws.Range("H126").Value = 4258
$ws.Range("I126").Value = 4301.3335
$ws.Range("J126").Value = 4180
$ws.Range("K126").Value = 12904.0005
$ws.Range("L126").Value = 12540
$ws.Range("M126").Value = -10434.0005
$ws.Range("N126").Value = -17480

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3641.739
$ws.Range("I7").Value = 3611.111
$ws.Range("J7").Value = 3661.4285
$ws.Range("K7").Value = 3611.111
$ws.Range("L7").Value = 3661.4285
$ws.Range("M7").Value = -3499.111
$ws.Range("N7").Value = -3885.4285
$ws.Range("H16").Value = 1300.0714
$ws.Range("I16").Value = 1447.1818
$ws.Range("J16").Value = 760.6667
$ws.Range("K16").Value = 1447.1818
$ws.Range("L16").Value = 760.6667
$ws.Range("M16").Value = -1277.1818
$ws.Range("N16").Value = -1100.6667
$ws.Range("H21").Value = 19666.666
$ws.Range("J21").Value = 19666.666
$ws.Range("L21").Value = 19666.666
$ws.Range("N21").Value = -20014.666
$ws.Range("H40").Value = 4884.875
$ws.Range("I40").Value = 6519.75
$ws.Range("J40").Value = 3250
$ws.Range("K40").Value = 6519.75
$ws.Range("L40").Value = 3250
$ws.Range("M40").Value = -6383.75
$ws.Range("N40").Value = -3522
$ws.Range("H76").Value = 9500
$ws.Range("I76").Value = 3000
$ws.Range("J76").Value = 10800
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 10800
$ws.Range("M76").Value = -2662
$ws.Range("N76").Value = -11476
$ws.Range("H79").Value = 9500
$ws.Range("I79").Value = 3000
$ws.Range("J79").Value = 10800
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 10800
$ws.Range("M79").Value = -1830
$ws.Range("N79").Value = -13140
$ws.Range("H104").Value = 20874.545
$ws.Range("J104").Value = 20874.545
$ws.Range("L104").Value = 20874.545
$ws.Range("N104").Value = -27862.545
$ws.Range("H122").Value = 7697027
$ws.Range("I122").Value = 3478.8572
$ws.Range("J122").Value = 16672833
$ws.Range("K122").Value = 10436.5716
$ws.Range("L122").Value = 50018499
$ws.Range("M122").Value = -7986.571599999999
$ws.Range("N122").Value = -50023399
$ws.Range("H126").Value = 3641.739
$ws.Range("I126").Value = 3611.111
$ws.Range("J126").Value = 3661.4285
$ws.Range("K126").Value = 10833.333
$ws.Range("L126").Value = 10984.2855
$ws.Range("M126").Value = -8363.332999999999
$ws.Range("N126").Value = -15924.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 47621856
$ws.Range("I122").Value = 76925040
$ws.Range("K122").Value = 230775120
$ws.Range("M122").Value = -230772670
